# Generate Report for Handback
# Refresh the localization-status report: update the "Status" text, bump the
# "Latest Handback DateTime" timestamps, clear the stale "Error Detail"
# message now that the handback is in sync, and widen a couple of columns
# that were too narrow for the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-28 12:49:59"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839

# ---- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-28 12:50:13"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839

Write-Output "Report regenerated for handback."
